$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.746.17"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.604.05"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.85"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.829.10"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.605.78"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.42"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.01"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  -4.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.70"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.36"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.286.15"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.22"
$ws.Range("E35").Value = "  +15.49%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("E37").Value = "  -4.74%  "
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.73"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.740.90"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.36"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.103"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.58"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +1.86%  "
